$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (row 1): shift "编号"/"名称" out, introduce
# "产品编号" and "版本号", and move the old "图号" header to column B.
$ws.Range("A1").Value = "产品编号"
$ws.Range("B1").Value = "图号"
$ws.Range("C1").Value = "版本号"

# Restore the active selection to a single cell A2 (was A6 / A2:A6).
$ws.Range("A2").Select()
